# Generate Report for Handoff
# - Bump the "Latest Handoff Datetime" timestamps for the zh-cn and de-de
#   handoff reports (rows 8-12 and 14, i.e. the current handoff batch).
# - The Overview sheet's "Latest HO Xliff Generate Date" column happened to
#   share the exact same text ("2016-08-15 22:20:36") as de-de's handoff
#   datetime for those same rows, so it moves in lockstep too.
# - Mark rows 8-12 and 14 (everything in the current handoff batch except the
#   "d1ba96db..." row) with Priority = "ht" on both language sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = 8, 9, 10, 11, 12, 14

foreach ($r in $rows) {
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"

    $zhcn.Range("H$r").Value = "2016-08-15 22:20:55"

    $dede.Range("H$r").Value = "2016-08-15 22:21:02"
    $overview.Range("G$r").Value = "2016-08-15 22:21:02"
}
